$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (2 through 138). Update it from 45178 (2023-09-09) to 45179 (2023-09-10),
# leaving the cell's existing number format/style untouched.
for ($row = 2; $row -le 138; $row++) {
    $ws.Cells.Item($row, 3).Value = 45179
}
